# 15.7.1.2 report: add a new "2023" data column (N) after the existing
# 2013-2022 year columns (D:M), carrying over the same header / border /
# value formatting as the last existing year column (M), and fill in the
# new year label + data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column M's formatting (and structure) for rows 3-5 into column N:
#   row 3 -> empty bottom-border cell (style "thin divider" row)
#   row 4 -> year header style
#   row 5 -> data value style
$ws.Range("M3:M5").Copy($ws.Range("N3:N5"))

# New year header and corresponding data value.
$ws.Range("N4").Value = 2023
$ws.Range("N5").Value = 553

# Row 3's height is nudged slightly as part of this edit.
$ws.Rows(3).RowHeight = 13.5
